# The edit re-orders two small blocks of species-observation rows on the
# "Artfynd" sheet: rows 10-13 cyclically shift up by one (row 10 takes row
# 11's data, ..., row 13 wraps around to take row 10's original data), and
# rows 18-22 cyclically shift down by one (row 18 takes row 22's original
# data, row 19 takes row 18's, ..., row 22 takes row 21's). Columns T/U/V/W
# (location), Y/AA (dates) and AD/AE/AG/AW/AX/AY (flags/reporter) are
# identical across all of these rows, so only the cells below actually
# change value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 10 (was row 11)
$ws.Range("A10").Value = 112501583
$ws.Range("B10").Value = 90830
$ws.Range("E10").Value = 2059
$ws.Range("F10").Value = 'Skrovlig taggsvamp'
$ws.Range("G10").Value = 'Hydnellum scabrosum'
$ws.Range("H10").Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("Q10").Value = 555598
$ws.Range("R10").Value = 7010805

# Row 11 (was row 12)
$ws.Range("A11").Value = 112501569
$ws.Range("B11").Value = 90808
$ws.Range("E11").Value = 4362
$ws.Range("F11").Value = 'Blå taggsvamp'
$ws.Range("G11").Value = 'Hydnellum caeruleum'
$ws.Range("H11").Value = '(Hornem.) P.Karst.'
$ws.Range("Q11").Value = 555402
$ws.Range("R11").Value = 7010819

# Row 12 (was row 13)
$ws.Range("A12").Value = 112501608
$ws.Range("B12").Value = 90814
$ws.Range("D12").Value = 'LC'
$ws.Range("E12").Value = 4364
$ws.Range("F12").Value = 'Dropptaggsvamp'
$ws.Range("G12").Value = 'Hydnellum ferrugineum'
$ws.Range("H12").Value = '(Fr.:Fr.) P. Karst.'
$ws.Range("Q12").Value = 555535
$ws.Range("R12").Value = 7011452

# Row 13 (was row 10)
$ws.Range("A13").Value = 112501614
$ws.Range("B13").Value = 90837
$ws.Range("D13").Value = 'NT'
$ws.Range("E13").Value = 5966
$ws.Range("F13").Value = 'Motaggsvamp'
$ws.Range("G13").Value = 'Sarcodon squamosus'
$ws.Range("H13").Value = '(Schaeff.) Quél.'
$ws.Range("Q13").Value = 555528
$ws.Range("R13").Value = 7011175

# Row 18 (was row 22)
$ws.Range("A18").Value = 112501617
$ws.Range("Q18").Value = 555555
$ws.Range("R18").Value = 7011181

# Row 19 (was row 18)
$ws.Range("A19").Value = 112501565
$ws.Range("B19").Value = 90808
$ws.Range("E19").Value = 4362
$ws.Range("F19").Value = 'Blå taggsvamp'
$ws.Range("G19").Value = 'Hydnellum caeruleum'
$ws.Range("H19").Value = '(Hornem.) P.Karst.'
$ws.Range("Q19").Value = 555424
$ws.Range("R19").Value = 7011077

# Row 20 (was row 19)
$ws.Range("A20").Value = 112501596
$ws.Range("B20").Value = 90837
$ws.Range("E20").Value = 5966
$ws.Range("F20").Value = 'Motaggsvamp'
$ws.Range("G20").Value = 'Sarcodon squamosus'
$ws.Range("H20").Value = '(Schaeff.) Quél.'
$ws.Range("Q20").Value = 555402
$ws.Range("R20").Value = 7011030
$ws.Range("AL20").Value = ""
$ws.Range("AO20").Value = ""

# Row 21 (was row 20)
$ws.Range("A21").Value = 112501575
$ws.Range("B21").Value = 77650
$ws.Range("E21").Value = 6425
$ws.Range("F21").Value = 'Garnlav'
$ws.Range("G21").Value = 'Alectoria sarmentosa'
$ws.Range("H21").Value = '(Ach.) Ach.'
$ws.Range("K21").Value = ""
$ws.Range("L21").Value = ""
$ws.Range("M21").Value = ""
$ws.Range("N21").Value = ""
$ws.Range("Q21").Value = 555531
$ws.Range("R21").Value = 7011013
$ws.Range("AC21").Value = ""
$ws.Range("AL21").Value = 'Tall'
$ws.Range("AO21").Value = 'Tall'

# Row 22 (was row 21)
$ws.Range("A22").Value = 112501584
$ws.Range("B22").Value = 56430
$ws.Range("E22").Value = 100109
$ws.Range("F22").Value = 'Tretåig hackspett'
$ws.Range("G22").Value = 'Picoides tridactylus'
$ws.Range("H22").Value = '(Linnaeus, 1758)'
$ws.Range("K22").Value = ""
$ws.Range("L22").Value = ""
$ws.Range("M22").Value = 'färska spår'
$ws.Range("N22").Value = ""
$ws.Range("Q22").Value = 555459
$ws.Range("R22").Value = 7011382
$ws.Range("AC22").Value = 'Hack'
